# NYPD 108th Precinct weekly CompStat report update.
# Rolls the reporting window forward one week (new crime data collected)
# and updates every statistic in the "Week to Date / 28 Day / Year to
# Date / 2 Year" crime-complaints table (rows 15-30 on the main sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Header text: "Volume 30 Number 9" -> "Volume 30 Number 10"
# ---------------------------------------------------------------------
$a8 = $ws.Range("A8").Text
$pos = $a8.LastIndexOf("9") + 1
$ws.Range("A8").Characters($pos, 1).Text = "10"

# ---------------------------------------------------------------------
# 2) Header text: report week "2/27/2023 ... 3/5/2023" -> "3/6/2023 ... 3/12/2023"
# ---------------------------------------------------------------------
$c9 = $ws.Range("C9").Text
$p1 = $c9.IndexOf("2/27/2023") + 1
$ws.Range("C9").Characters($p1, 9).Text = "3/6/2023"

$c9b = $ws.Range("C9").Text
$p2 = $c9b.IndexOf("3/5/2023") + 1
$ws.Range("C9").Characters($p2, 8).Text = "3/12/2023"

# ---------------------------------------------------------------------
# Helper: set a plain numeric value, forcing the right display format
# when the cell's format needs to change too (string->number swaps).
# ---------------------------------------------------------------------
function Set-Num($addr, $value, $fmt) {
    $ws.Range($addr).Value = $value
    if ($fmt) {
        $ws.Range($addr).NumberFormat = $fmt
    }
}

$FMT_COUNT = "#,##0"
$FMT_PCT   = "#,##0.0;`"-`"#,##0.0"

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
Set-Num "C15" 2
Set-Num "F15" 4
Set-Num "I15" 11
Set-Num "K15" 266.666666666667
Set-Num "L15" 1000
Set-Num "M15" 450
Set-Num "N15" 1000

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
Set-Num "D16" 4
Set-Num "E16" 25
Set-Num "G16" 20
Set-Num "H16" 5
Set-Num "I16" 51
Set-Num "J16" 37
Set-Num "K16" 37.837837837837
Set-Num "L16" 200
Set-Num "M16" 18.604651162790
Set-Num "N16" -79.518072289156

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
Set-Num "C17" 4
Set-Num "D17" 1
Set-Num "E17" 300
Set-Num "G17" 12
Set-Num "H17" -16.666666666666
Set-Num "I17" 34
Set-Num "J17" 34
Set-Num "K17" 0
Set-Num "L17" 9.677419354838
Set-Num "M17" 30.769230769230
Set-Num "N17" -29.166666666666

# ---------------------------------------------------------------------
# Row 18 - Burglary (D18/E18 flip from text placeholders to real numbers)
# ---------------------------------------------------------------------
Set-Num "C18" 2
Set-Num "D18" 1 $FMT_COUNT
Set-Num "E18" 100 $FMT_PCT
Set-Num "F18" 16
Set-Num "G18" 8
Set-Num "H18" 100
Set-Num "I18" 39
Set-Num "J18" 33
Set-Num "K18" 18.181818181818
Set-Num "L18" 21.875
Set-Num "M18" -18.75
Set-Num "N18" -86.315789473684

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
Set-Num "D19" 14
Set-Num "E19" 14.285714285714
Set-Num "F19" 53
Set-Num "G19" 47
Set-Num "H19" 12.765957446808
Set-Num "I19" 138
Set-Num "J19" 138
Set-Num "L19" 91.666666666666
Set-Num "M19" 84
Set-Num "N19" -6.756756756756

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
Set-Num "C20" 5
Set-Num "D20" 8
Set-Num "E20" -37.5
Set-Num "F20" 15
Set-Num "G20" 12
Set-Num "H20" 25
Set-Num "I20" 37
Set-Num "J20" 34
Set-Num "K20" 8.823529411764
Set-Num "L20" 27.586206896551
Set-Num "M20" -19.565217391304
Set-Num "N20" -90.389610389610

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
Set-Num "C21" 34
Set-Num "D21" 28
Set-Num "E21" 21.428571428571
Set-Num "F21" 119
Set-Num "G21" 99
Set-Num "H21" 20.202020202020
Set-Num "I21" 310
Set-Num "J21" 279
Set-Num "K21" 11.111111111111
Set-Num "L21" 70.329670329670
Set-Num "M21" 28.099173553719
Set-Num "N21" -72.296693476318

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
Set-Num "D22" 2
Set-Num "E22" -50
Set-Num "F22" 6
Set-Num "G22" 9
Set-Num "H22" -33.333333333333
Set-Num "I22" 14
Set-Num "J22" 18
Set-Num "K22" -22.222222222222
Set-Num "L22" 600
Set-Num "M22" 55.555555555555

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
Set-Num "C24" 40
Set-Num "D24" 17
Set-Num "E24" 135.294117647059
Set-Num "F24" 162
Set-Num "G24" 92
Set-Num "H24" 76.086956521739
Set-Num "I24" 331
Set-Num "J24" 212
Set-Num "K24" 56.132075471698
Set-Num "L24" 64.676616915422
Set-Num "M24" 138.129496402878

# ---------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------
Set-Num "C25" 6
Set-Num "E25" -50
Set-Num "F25" 33
Set-Num "G25" 51
Set-Num "H25" -35.294117647058
Set-Num "I25" 103
Set-Num "J25" 114
Set-Num "K25" -9.649122807017
Set-Num "L25" 66.129032258064
Set-Num "M25" 3

# ---------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------
Set-Num "C26" 2
Set-Num "F26" 4
Set-Num "I26" 13
Set-Num "K26" 333.333333333333
Set-Num "L26" 550

# ---------------------------------------------------------------------
# Row 27 - Other Sex Crimes (C27 flips from text placeholder to a number)
# ---------------------------------------------------------------------
Set-Num "C27" 1 $FMT_COUNT
Set-Num "D27" 1
Set-Num "E27" 0
Set-Num "F27" 4
Set-Num "G27" 7
Set-Num "H27" -42.857142857142
Set-Num "I27" 18
Set-Num "J27" 13
Set-Num "K27" 38.461538461538
Set-Num "L27" 125

# ---------------------------------------------------------------------
# Row 30 - Hate Crimes (C30/D30/E30 flip from real numbers to the
# workbook's "no data" text placeholders: "0" / "***.*")
# ---------------------------------------------------------------------
function Set-TextPlaceholder($addr, $text) {
    $ws.Range($addr).Style = "Normal"
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Font.Name = "Andale WT"
    $ws.Range($addr).Font.Size = 10
    $ws.Range($addr).HorizontalAlignment = -4152
    $ws.Range($addr).VerticalAlignment = -4108
}

Set-TextPlaceholder "C30" "0"
Set-TextPlaceholder "D30" "0"
Set-TextPlaceholder "E30" "***.*"
